$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: update title (D26)
$ws.Range("D26").Value = "2021 인공지능 경진대회 참가기"

# Row 28: update title (D28) and link (E28)
$ws.Range("D28").Value = "Let's do MuJoCo - 2.2 Mujoco XML File(MJCF) 실행"
$ws.Range("E28").Value = "https://ropiens.tistory.com/171"

# Row 51: update title (D51) and link (E51)
$ws.Range("D51").Value = "[MariaDB] MariaDB 서버 타임존 변경하기"
$ws.Range("E51").Value = "https://bskyvision.com/1236"
